$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# The document (before the <w:tbl>) has 7 paragraphs:
#   1) "Caso de uso Ver Tablón de anuncios "           -> becomes bold title
#   2) "" (empty)                                       -> becomes "Caso de uso: Ver ..." + "Pasos:"
#   3) numbered "Hacer click en el botón de “Home”"     -> unchanged
#   4) numbered "Se muestran los mensajes..."           -> unchanged
#   5) "" (empty, ind left=720)                         -> unchanged, but new "Precondiciones:" p. inserted before it
#   6) "" (empty)                                        -> unchanged
#   7) "" (empty)                                        -> becomes bold "Caso de Prueba:"
#
# We edit from the bottom paragraph upward so earlier indices stay valid.
# ---------------------------------------------------------------------------

# --- Paragraph 7 (last empty paragraph before the table): "Caso de Prueba:" ---
$p7 = $d.Paragraphs.Item(7)
$xmlCasoPrueba = @"
<w:p $wNs>
  <w:pPr>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
  </w:pPr>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t xml:space="preserve">Caso </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t>de Prueba</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t>:</w:t>
  </w:r>
</w:p>
"@
$p7.Range.InsertXML($xmlCasoPrueba)

# --- Paragraph 5 (the "ind left=720" empty one): insert a new "Precondiciones:" paragraph
#     immediately before it, while preserving paragraph 5 itself unchanged. Using a
#     collapsed range at Paragraph 5's start merges into paragraph 5 in this host, so
#     instead we replace paragraph 5's whole range with [new paragraph, original paragraph]. ---
$p5 = $d.Paragraphs.Item(5)
$xmlPrecondPlusOriginal = @"
<w:p $wNs>
  <w:r>
    <w:t>Precondiciones:</w:t>
  </w:r>
  <w:r>
    <w:br/>
    <w:t xml:space="preserve">     1) Estar registrado y </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>logueado</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t xml:space="preserve"> en el sistema. </w:t>
  </w:r>
</w:p>
<w:p $wNs>
  <w:pPr>
    <w:ind w:left="720"/>
  </w:pPr>
</w:p>
"@
$p5.Range.InsertXML($xmlPrecondPlusOriginal)

# --- Paragraph 1: "Caso de uso Ver Tablón de anuncios " -> bold "Procedimiento de Pruebas:" ---
$p1 = $d.Paragraphs.Item(1)
$xmlTitle = @"
<w:p $wNs>
  <w:pPr>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t>Procedimiento de Pruebas:</w:t>
  </w:r>
</w:p>
"@
$p1.Range.InsertXML($xmlTitle)

# --- Paragraph 2 (empty) -> "Caso de uso: Ver Tablón de anuncios " + "Pasos:" (two paragraphs) ---
$p2 = $d.Paragraphs.Item(2)
$xmlCasoUsoYPasos = @"
<w:p $wNs>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t>Caso de uso</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t>:</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> Ver Tablón de anuncios </w:t>
  </w:r>
</w:p>
<w:p $wNs>
  <w:pPr>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t>Pasos:</w:t>
  </w:r>
</w:p>
"@
$p2.Range.InsertXML($xmlCasoUsoYPasos)

Write-Output "Edit complete."
